$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column cells being touched stay as text (matches original
# inline-string storage) instead of being auto-converted to numbers by Excel.
$priceCells = @("D2","D3","D5","D8","D12","D13","D16","D17","D18","D19","D23","D25","D26","D34","D45","D46","D47","D49","D51")
foreach ($cellRef in $priceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.426.69"
$ws.Range("E2").Value = "  +3.41%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.591.33"
$ws.Range("E3").Value = "  +1.47%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.96%  "

# Row 5 - BNB
$ws.Range("D5").Value = "213.15"
$ws.Range("E5").Value = "  +0.92%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.12%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.95%  "

# Row 8 - Solana
$ws.Range("D8").Value = "24.40"
$ws.Range("E8").Value = "  +7.62%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +0.49%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.67%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.818.03"
$ws.Range("E12").Value = "  +1.49%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.596.84"
$ws.Range("E13").Value = "  +2.10%  "

# Row 14 - Polygon
$ws.Range("E14").Value = "  +2.15%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  -0.15%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "28.451.48"
$ws.Range("E16").Value = "  +3.57%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "63.13"
$ws.Range("E17").Value = "  +1.20%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "229.67"
$ws.Range("E18").Value = "  +1.57%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0706"
$ws.Range("E19").Value = "  +0.25%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  -0.35%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.87%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -1.13%  "

# Row 23 - Avalanche
$ws.Range("D23").Value = "9.33"
$ws.Range("E23").Value = "  -0.66%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  +1.14%  "

# Row 25 - Monero
$ws.Range("D25").Value = "151.68"
$ws.Range("E25").Value = "  +1.14%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "15.23"
$ws.Range("E26").Value = "  +0.65%  "

# Row 27 - Cosmos
$ws.Range("E27").Value = "  -0.74%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -0.65%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.93%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.82%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +0.42%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.38%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +0.88%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.401.56"
$ws.Range("E34").Value = "  -3.56%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -0.78%  "

# Row 36 - TrustWalletToken
$ws.Range("E36").Value = "  -9.58%  "

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  +1.02%  "

# Row 38 - MXToken
$ws.Range("E38").Value = "  +8.76%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  -0.57%  "

# Row 40 - ImmutableX
$ws.Range("E40").Value = "  +0.35%  "

# Row 41 - ARBITRUM
$ws.Range("E41").Value = "  -0.03%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.95%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  -2.22%  "

# Row 44 - RenderToken
$ws.Range("E44").Value = "  +0.29%  "

# Row 45 - WEMIXToken
$ws.Range("D45").Value = "0.982"
$ws.Range("E45").Value = "  +0.94%  "

# Row 46 - Aave
$ws.Range("D46").Value = "63.19"
$ws.Range("E46").Value = "  -1.71%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.727.44"
$ws.Range("E47").Value = "  +1.46%  "

# Row 48 - mCoin
$ws.Range("E48").Value = "  +1.73%  "

# Row 49 - Quant
$ws.Range("D49").Value = "87.21"
$ws.Range("E49").Value = "  +0.56%  "

# Row 50 - BabyDogeCoin
$ws.Range("E50").Value = "  +1.08%  "

# Row 51 - Cronos
$ws.Range("D51").Value = "0.0521"
$ws.Range("E51").Value = "  -0.78%  "
